# Add a new "Kestrekel_AC" animal companion worksheet.
#
# The new sheet is structurally identical to the existing "Kivit_AC" sheet
# (same columns/widths, same layout, same formulas) but with a handful of
# the "Feat" cells in column E swapped out for the Kes'trekel's own feats.
# It is inserted as the last tab (after "Z'tal_AC") and becomes the active
# sheet, so it steals `tabSelected` away from "Kivit_AC".

$wb = $excel.ActiveWorkbook

$source = $wb.Worksheets.Item("Kivit_AC")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# Copy Kivit_AC to the end of the workbook; Excel makes the new copy the
# active sheet, which also clears tabSelected on the source sheet.
$source.Copy($null, $lastSheet)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Kestrekel_AC"

# Kes'trekel-specific feat progression (differs from Kivit_AC in these cells
# only; two of the values are brand-new feats not used anywhere else).
$newSheet.Range("E6").Value = "Spring Attack*"
$newSheet.Range("E7").Value = "WF(Creature)"
$newSheet.Range("E9").Value = "Dodge"
$newSheet.Range("E16").Value = "Improved Critical"
$newSheet.Range("E19").Value = "Epic Prowess"

$newSheet.Range("F3").Select()
